$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace header labels with numeric column indices (0-based), keep existing style
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# Row 2: move the header labels here (this row previously held "Grade 2 Titanium" in A2 only)
$ws.Range("A2").Value = "Lg."
$ws.Range("B2").Value = "Threading"
$ws.Range("C2").Value = "HeadDia."
$ws.Range("D2").Value = "Head Ht."
$ws.Range("E2").Value = "DriveSize"
$ws.Range("F2").Value = "TensileStrength, psi"
$ws.Range("G2").Value = "Specifications Met"
$ws.Range("H2").Value = "Pkg.Qty."
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Pkg."
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""

# Column L, rows 3-27: fill with the material that used to sit in A2
for ($r = 3; $r -le 27; $r++) {
    $ws.Cells.Item($r, 12).Value = "Grade 2 Titanium"
}
